$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy A/E column formatting (bold/centered id style, date style) down to the new rows 219-226
$ws.Range("A213").Copy() | Out-Null
$ws.Range("A219:A226").PasteSpecial(-4122) | Out-Null
$ws.Range("E213").Copy() | Out-Null
$ws.Range("E219:E226").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 213
$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 211
$arr[0,1] = 7973582
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45380.375
$arr[0,5] = "Pirin Blagoevgrad"
$arr[0,6] = "Etar 1924 Veliko Tarnovo"
$arr[0,7] = 2
$arr[0,8] = 1
$arr[0,9] = "H"
$arr[0,10] = 1.909
$arr[0,11] = 3.3
$arr[0,12] = 4.2
$arr[0,13] = 2
$arr[0,14] = 3.3
$arr[0,15] = 4
$arr[0,16] = -0.5
$arr[0,17] = 2.025
$arr[0,18] = 1.825
$arr[0,19] = 2
$arr[0,20] = 1.975
$arr[0,21] = 1.875
$arr[0,22] = 1
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 1.025
$arr[0,26] = -1
$arr[0,27] = 0.9750000000000001
$arr[0,28] = -1
$ws.Range("A213:AC213").Value2 = $arr

# Row 214
$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 212
$arr[0,1] = 7973583
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45380.47916666666
$arr[0,5] = "Arda Kardzhali"
$arr[0,6] = "Krumovgrad"
$arr[0,7] = 1
$arr[0,8] = 1
$arr[0,9] = "D"
$arr[0,10] = 1.833
$arr[0,11] = 3.4
$arr[0,12] = 4.5
$arr[0,13] = 2.1
$arr[0,14] = 3.1
$arr[0,15] = 4
$arr[0,16] = -0.25
$arr[0,17] = 1.8
$arr[0,18] = 2.05
$arr[0,19] = 2
$arr[0,20] = 1.875
$arr[0,21] = 1.975
$arr[0,22] = -1
$arr[0,23] = 2.1
$arr[0,24] = -1
$arr[0,25] = -0.5
$arr[0,26] = 0.5249999999999999
$arr[0,27] = 0
$arr[0,28] = -0
$ws.Range("A214:AC214").Value2 = $arr

# Row 215
$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 213
$arr[0,1] = 7973584
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45380.58333333334
$arr[0,5] = "Levski Sofia"
$arr[0,6] = "Slavia Sofia"
$arr[0,7] = 2
$arr[0,8] = 0
$arr[0,9] = "H"
$arr[0,10] = 1.444
$arr[0,11] = 4.2
$arr[0,12] = 7.5
$arr[0,13] = 1.6
$arr[0,14] = 3.75
$arr[0,15] = 6.5
$arr[0,16] = -0.75
$arr[0,17] = 1.8
$arr[0,18] = 2.05
$arr[0,19] = 2
$arr[0,20] = 1.925
$arr[0,21] = 1.925
$arr[0,22] = 0.6000000000000001
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 0.8
$arr[0,26] = -1
$arr[0,27] = 0
$arr[0,28] = -0
$ws.Range("A215:AC215").Value2 = $arr

# Row 216
$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 214
$arr[0,1] = 7973585
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45381.375
$arr[0,5] = "Botev Vratsa"
$arr[0,6] = "FC Hebar Pazardzhik"
$arr[0,7] = 3
$arr[0,8] = 2
$arr[0,9] = "H"
$arr[0,10] = 2.4
$arr[0,11] = 3.2
$arr[0,12] = 3
$arr[0,13] = 2.7
$arr[0,14] = 3.1
$arr[0,15] = 2.8
$arr[0,16] = 0
$arr[0,17] = 1.875
$arr[0,18] = 1.975
$arr[0,19] = 2
$arr[0,20] = 1.875
$arr[0,21] = 1.975
$arr[0,22] = 1.7
$arr[0,23] = -1
$arr[0,24] = -1
$arr[0,25] = 0.875
$arr[0,26] = -1
$arr[0,27] = 0.875
$arr[0,28] = -1
$ws.Range("A216:AC216").Value2 = $arr

# Row 217
$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 215
$arr[0,1] = 7973586
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45381.47916666666
$arr[0,5] = "Lokomotiv Plovdiv"
$arr[0,6] = "Ludogorets Razgrad"
$arr[0,7] = 1
$arr[0,8] = 2
$arr[0,9] = "A"
$arr[0,10] = 5.25
$arr[0,11] = 3.75
$arr[0,12] = 1.65
$arr[0,13] = 8.5
$arr[0,14] = 4.333
$arr[0,15] = 1.4
$arr[0,16] = 1.25
$arr[0,17] = 1.925
$arr[0,18] = 1.925
$arr[0,19] = 2.75
$arr[0,20] = 2.025
$arr[0,21] = 1.825
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.3999999999999999
$arr[0,25] = 0.4625
$arr[0,26] = -0.5
$arr[0,27] = 0.5125
$arr[0,28] = -0.5
$ws.Range("A217:AC217").Value2 = $arr

# Row 218
$arr = New-Object 'object[,]' 1,29
$arr[0,0] = 216
$arr[0,1] = 6978434
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45381.58333333334
$arr[0,5] = "Lokomotiv 1929 Sofia"
$arr[0,6] = "CSKA Sofia"
$arr[0,7] = 0
$arr[0,8] = 1
$arr[0,9] = "A"
$arr[0,10] = 7.5
$arr[0,11] = 5
$arr[0,12] = 1.363
$arr[0,13] = 12
$arr[0,14] = 5.75
$arr[0,15] = 1.25
$arr[0,16] = 1.5
$arr[0,17] = 2.025
$arr[0,18] = 1.825
$arr[0,19] = 2.5
$arr[0,20] = 1.95
$arr[0,21] = 1.9
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.25
$arr[0,25] = 1.025
$arr[0,26] = -1
$arr[0,27] = -1
$arr[0,28] = 0.8999999999999999
$ws.Range("A218:AC218").Value2 = $arr

# Row 219
$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 217
$arr[0,1] = 6978445
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45384.45833333334
$arr[0,5] = "Etar 1924 Veliko Tarnovo"
$arr[0,6] = "Arda Kardzhali"
$arr[0,10] = 4.5
$arr[0,11] = 3.4
$arr[0,12] = 1.833
$arr[0,13] = 4.333
$arr[0,14] = 3.4
$arr[0,15] = 1.909
$arr[0,16] = 0.5
$arr[0,17] = 1.925
$arr[0,18] = 1.925
$arr[0,19] = 2.25
$arr[0,20] = 1.9
$arr[0,21] = 1.95
$arr[0,22] = 0
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 0
$arr[0,26] = 0
$ws.Range("A219:AA219").Value2 = $arr

# Row 220
$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 218
$arr[0,1] = 6978444
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45384.5625
$arr[0,5] = "Slavia Sofia"
$arr[0,6] = "Pirin Blagoevgrad"
$arr[0,10] = 1.571
$arr[0,11] = 3.8
$arr[0,12] = 6
$arr[0,13] = 1.5
$arr[0,14] = 4
$arr[0,15] = 7.5
$arr[0,16] = -1
$arr[0,17] = 1.875
$arr[0,18] = 1.975
$arr[0,19] = 2.25
$arr[0,20] = 1.9
$arr[0,21] = 1.95
$arr[0,22] = 0
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 0
$arr[0,26] = 0
$ws.Range("A220:AA220").Value2 = $arr

# Row 221
$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 219
$arr[0,1] = 6978389
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45385.35416666666
$arr[0,5] = "Krumovgrad"
$arr[0,6] = "FC Hebar Pazardzhik"
$arr[0,10] = 1.909
$arr[0,11] = 3.3
$arr[0,12] = 4.2
$arr[0,13] = 2
$arr[0,14] = 3.25
$arr[0,15] = 4
$arr[0,16] = -0.5
$arr[0,17] = 2.025
$arr[0,18] = 1.825
$arr[0,19] = 2
$arr[0,20] = 1.875
$arr[0,21] = 1.975
$arr[0,22] = 0
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 0
$arr[0,26] = 0
$ws.Range("A221:AA221").Value2 = $arr

# Row 222
$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 220
$arr[0,1] = 6978442
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45385.45833333334
$arr[0,5] = "CSKA Sofia"
$arr[0,6] = "Lokomotiv Plovdiv"
$arr[0,10] = 1.45
$arr[0,11] = 4
$arr[0,12] = 8
$arr[0,13] = 1.333
$arr[0,14] = 4.5
$arr[0,15] = 12
$arr[0,16] = -1.25
$arr[0,17] = 1.825
$arr[0,18] = 2.025
$arr[0,19] = 2.25
$arr[0,20] = 1.8
$arr[0,21] = 2.05
$arr[0,22] = 0
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 0
$arr[0,26] = 0
$ws.Range("A222:AA222").Value2 = $arr

# Row 223
$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 221
$arr[0,1] = 6978443
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45385.5625
$arr[0,5] = "Ludogorets Razgrad"
$arr[0,6] = "Levski Sofia"
$arr[0,10] = 1.571
$arr[0,11] = 3.8
$arr[0,12] = 6
$arr[0,13] = 1.5
$arr[0,14] = 4
$arr[0,15] = 7
$arr[0,16] = -1
$arr[0,17] = 1.875
$arr[0,18] = 1.975
$arr[0,19] = 2.25
$arr[0,20] = 1.825
$arr[0,21] = 2.025
$arr[0,22] = 0
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 0
$arr[0,26] = 0
$ws.Range("A223:AA223").Value2 = $arr

# Row 224
$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 222
$arr[0,1] = 6978439
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45386.35416666666
$arr[0,5] = "CSKA 1948 Sofia"
$arr[0,6] = "Botev Vratsa"
$arr[0,10] = 1.333
$arr[0,11] = 5.5
$arr[0,12] = 7.5
$arr[0,13] = 1.444
$arr[0,14] = 4.75
$arr[0,15] = 5.5
$arr[0,16] = -1.25
$arr[0,17] = 2.05
$arr[0,18] = 1.8
$arr[0,19] = 2.5
$arr[0,20] = 1.95
$arr[0,21] = 1.9
$arr[0,22] = 0
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 0
$arr[0,26] = 0
$ws.Range("A224:AA224").Value2 = $arr

# Row 225
$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 223
$arr[0,1] = 6978440
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45386.45833333334
$arr[0,5] = "Beroe"
$arr[0,6] = "Botev Plovdiv"
$arr[0,10] = 3.5
$arr[0,11] = 3.2
$arr[0,12] = 2.15
$arr[0,13] = 4.333
$arr[0,14] = 3.3
$arr[0,15] = 1.85
$arr[0,16] = 0.5
$arr[0,17] = 1.925
$arr[0,18] = 1.925
$arr[0,19] = 2.5
$arr[0,20] = 2.05
$arr[0,21] = 1.8
$arr[0,22] = 0
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 0
$arr[0,26] = 0
$ws.Range("A225:AA225").Value2 = $arr

# Row 226
$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 224
$arr[0,1] = 6978441
$arr[0,2] = "Bulgaria First League"
$arr[0,3] = "Bulgaria First League"
$arr[0,4] = 45386.5625
$arr[0,5] = "Cherno More Varna"
$arr[0,6] = "Lokomotiv 1929 Sofia"
$arr[0,10] = 1.363
$arr[0,11] = 4.75
$arr[0,12] = 8.5
$arr[0,13] = 1.333
$arr[0,14] = 4.333
$arr[0,15] = 12
$arr[0,16] = -1.5
$arr[0,17] = 2.05
$arr[0,18] = 1.8
$arr[0,19] = 2.25
$arr[0,20] = 1.825
$arr[0,21] = 2.025
$arr[0,22] = 0
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 0
$arr[0,26] = 0
$ws.Range("A226:AA226").Value2 = $arr
